$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Adez Warm Shelf Share" to "Adez Shelf Share"
$ws.Range("A44").Value = "Adez Shelf Share"

# Remove the "Adez Cold Shelf Share" row entirely (row 45)
$ws.Rows("45").Delete()

# Move selection to reflect the edit location
$ws.Range("B44").Select()
